$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2022-04-17"

# Update the 2022 column header label (column I, row 1)
$ws.Range("I1").Value = "2022 (through 04-17)"

# Update April 2022 total (row 5 = April)
$ws.Range("I5").Value = 74

# Update December 2021 value (row 13 = December)
$ws.Range("H13").Value = 204

# Update 2021 yearly total (row 14)
$ws.Range("H14").Value = 1851

# Update 2022 yearly total (row 14)
$ws.Range("I14").Value = 509
